# Update "faturamento_anual" data for year 2025 (row 9) on the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B9").Value = 3154975.38
$ws.Range("C9").Value = 494771.99
$ws.Range("D9").Value = 3649747.37
$ws.Range("E9").Value = 13.5563352704053
$ws.Range("F9").Value = 86.44366472959469
$ws.Range("G9").Value = -52.1827542133771
$ws.Range("H9").Value = -43.02550556915311
$ws.Range("I9").Value = 31392
$ws.Range("J9").Value = 1333
$ws.Range("K9").Value = 32725
$ws.Range("L9").Value = 22608
$ws.Range("M9").Value = 161.4361009377212
$ws.Range("N9").Value = 10.21571419700311
